{"js": "// Mark two finished TODO list items as done (strike-through) and add a new,\n// still-open TODO item (\"Einheiten frei konfigurierbar machen\") right before\n// the trailing bookmark paragraph, followed by a fresh blank list paragraph.\n\nconst HOCH_TEXT = \"Hoch- und Runterschieben von \u00dcbungen im Training implementieren\";\nconst AUSBLENDEN_TEXT = \"Ausblenden Trianingssettings nach Start\";\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\nlet hochParagraph = null;\nlet ausblendenParagraph = null;\nlet bookmarkParagraph = null;\n\nconst items = paragraphs.items;\nfor (let i = 0; i < items.length; i++) {\n  const text = items[i].text;\n  if (text === HOCH_TEXT) {\n    hochParagraph = items[i];\n  } else if (text === AUSBLENDEN_TEXT) {\n    ausblendenParagraph = items[i];\n    // The empty paragraph immediately following \"Ausblenden...\" carries\n    // the _GoBack bookmark; that is the paragraph we expand below.\n    bookmarkParagraph = items[i + 1];\n  }\n}\n\nif (!hochParagraph) {\n  throw new Error(\"Could not find paragraph: \" + HOCH_TEXT);\n}\nif (!ausblendenParagraph || !bookmarkParagraph) {\n  throw new Error(\"Could not find paragraph: \" + AUSBLENDEN_TEXT + \" (or its bookmark sibling)\");\n}\n\n// 1) Strike through the two completed items. Setting font.strikeThrough on\n//    the paragraph applies <w:strike/> to both the run(s) and the\n//    paragraph mark, matching the target markup.\nhochParagraph.font.strikeThrough = true;\nausblendenParagraph.font.strikeThrough = true;\n\n// 2) Replace the (currently empty) bookmark paragraph with its new\n//    content: a new, still-open TODO item (its text is NOT struck\n//    through, only the paragraph mark is) that keeps the _GoBack\n//    bookmark, followed by a brand-new empty list paragraph.\nfunction wrapAsPackageOoxml(bodyXml) {\n  return (\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"512\">' +\n    '<pkg:xmlData><Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n    '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' +\n    \"</Relationships></pkg:xmlData></pkg:part>\" +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' +\n    bodyXml +\n    '<w:sectPr><w:pgSz w:w=\"12240\" w:h=\"15840\"/></w:sectPr>' +\n    \"</w:body></w:document></pkg:xmlData></pkg:part>\" +\n    \"</pkg:package>\"\n  );\n}\n\nconst newParagraphsXml = wrapAsPackageOoxml(\n  \"<w:p>\" +\n    \"<w:pPr>\" +\n    '<w:pStyle w:val=\"ListParagraph\"/>' +\n    '<w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr>' +\n    \"<w:rPr><w:strike/></w:rPr>\" +\n    \"</w:pPr>\" +\n    \"<w:r><w:t>Einheiten frei konfigurierbar machen</w:t></w:r>\" +\n    '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>' +\n    \"</w:p>\" +\n    '<w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/></w:pPr></w:p>'\n);\n\nbookmarkParagraph.insertOoxml(newParagraphsXml, Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# Mark two finished TODO list items as done (strike-through) and add a new,\n# still-open TODO item (\"Einheiten frei konfigurierbar machen\") right before\n# the trailing bookmark paragraph, followed by a fresh blank list paragraph.\n\n$d = $word.ActiveDocument\n\n$HOCH_TEXT = \"Hoch- und Runterschieben von \u00dcbungen im Training implementieren\"\n$AUSBLENDEN_TEXT = \"Ausblenden Trianingssettings nach Start\"\n\n$hochIndex = -1\n$ausblendenIndex = -1\n\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    # Paragraph.Range.Text includes the trailing paragraph mark (chr 13);\n    # trim it (and the rare chr 7 cell-mark) before comparing.\n    $text = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)\n    if ($text -eq $HOCH_TEXT) {\n        $hochIndex = $i\n    } elseif ($text -eq $AUSBLENDEN_TEXT) {\n        $ausblendenIndex = $i\n    }\n}\n\nif ($hochIndex -eq -1) {\n    throw \"Could not find paragraph: $HOCH_TEXT\"\n}\nif ($ausblendenIndex -eq -1) {\n    throw \"Could not find paragraph: $AUSBLENDEN_TEXT\"\n}\n\n# The empty paragraph immediately following \"Ausblenden...\" carries the\n# _GoBack bookmark; that is the paragraph we expand below.\n$bookmarkIndex = $ausblendenIndex + 1\n\n# 1) Strike through the two completed items. Setting StrikeThrough on the\n#    paragraph's Range applies <w:strike/> to both the run(s) and the\n#    paragraph mark, matching the target markup.\n$d.Paragraphs.Item($hochIndex).Range.Font.StrikeThrough = 1\n$d.Paragraphs.Item($ausblendenIndex).Range.Font.StrikeThrough = 1\n\n# 2) Replace the (currently empty) bookmark paragraph with its new content:\n#    a new, still-open TODO item (its text is NOT struck through, only the\n#    paragraph mark is) that keeps the _GoBack bookmark, followed by a\n#    brand-new empty list paragraph.\n$newParagraphsXml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"512\">' +\n    '<pkg:xmlData><Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n    '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' +\n    '</Relationships></pkg:xmlData></pkg:part>' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' +\n    '<w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr><w:rPr><w:strike/></w:rPr></w:pPr>' +\n    '<w:r><w:t>Einheiten frei konfigurierbar machen</w:t></w:r>' +\n    '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>' +\n    '</w:p>' +\n    '<w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/></w:pPr></w:p>' +\n    '<w:sectPr><w:pgSz w:w=\"12240\" w:h=\"15840\"/></w:sectPr>' +\n    '</w:body></w:document></pkg:xmlData></pkg:part>' +\n    '</pkg:package>'\n\n$d.Paragraphs.Item($bookmarkIndex).Range.InsertXML($newParagraphsXml)\n"}
